$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("oral")

# --- Copy the formatting of the last existing data row (row 24) onto the
# --- three new rows (25-27) so number formats / alignment match the rest
# --- of the table (text-as-number columns, right/left alignment, date
# --- format on column L, etc.) before any values are written.
$ws.Range("A24:L24").Copy() | Out-Null
$ws.Range("A25:L25").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:L26").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:L27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Shared "conference" / "location" / "day" values for the new presentations.
$ws.Range("C25").Value = "6th CYSTINET Working Group Meeting"
$ws.Range("C26").Value = "6th CYSTINET Working Group Meeting"
$ws.Range("G25").Value = "Ljubljana, Slovenia"
$ws.Range("G26").Value = "Ljubljana, Slovenia"
$ws.Range("G27").Value = "Ljubljana, Slovenia"

# Row 25 - Hobbs et al.
$ws.Range("A25").Value = "Assessment of the computer-based Taenia solium educational program 'The Vicious Worm' on knowledge update in primary school children in Katete district in Eastern Zambia"
$ws.Range("B25").Value = "Hobbs, E C; Mwape, K E; Zulu, G; Mambwe, M; Chembensofu, M; Phiri, I K; Masuku, M; Berkvens, D; Bottieau, E; Devleesschauwer, Brecht; Speybroeck, Niko; Colston, A; Willingham, A L; Dorny, Pierre; Gabriël, Sarah"

# Row 26 - Dupont et al.
$ws.Range("A26").Value = "Risk factor analysis in patients with neurocysticercosis associated epilepsy in northern Uganda"
$ws.Range("B26").Value = "Dupont, Fabian; Devleesschauwer, Brecht; Kaducu, Joyce; Lauseker, M; Schmidt, Veronika; Ovuga, E; Winkler, Andrea"

# Row 27 - Bouwknegt et al.
$ws.Range("A27").Value = "Ranking foodborne parasites in Europe using multicriteria decision analyses"
$ws.Range("B27").Value = "Bouwknegt, Martijn; Graham, Heather; Devleesschauwer, Brecht; Robertson, Lucy; van der Giessen, Joke"
$ws.Range("C27").Value = "Joint EURO-FBP and CYSTINET Meeting"

# "day" column
$ws.Range("F25").Value = "27"
$ws.Range("F26").Value = "27"
$ws.Range("F27").Value = "28"

# Row 25 remaining fields
$ws.Range("D25").Value = 2016
$ws.Range("E25").Value = 9
$ws.Range("H25").Value = "NA"
$ws.Range("I25").Value = "NA"
$ws.Range("J25").Value = "NA"
$ws.Range("K25").Value = "NA"
$ws.Range("L25").Value = 42640

# Row 26 remaining fields
$ws.Range("D26").Value = 2016
$ws.Range("E26").Value = 9
$ws.Range("H26").Value = "NA"
$ws.Range("I26").Value = "NA"
$ws.Range("J26").Value = "NA"
$ws.Range("K26").Value = "NA"
$ws.Range("L26").Value = 42640

# Row 27 remaining fields
$ws.Range("D27").Value = 2016
$ws.Range("E27").Value = 9
$ws.Range("H27").Value = "NA"
$ws.Range("I27").Value = "NA"
$ws.Range("J27").Value = "NA"
$ws.Range("K27").Value = "NA"
$ws.Range("L27").Value = 42641

# --- Grow the "Tabel1" table (and its autofilter) so it covers the new rows.
$lo = $ws.ListObjects.Item("Tabel1")
$lo.Resize($ws.Range("A1:L27"))

# --- Zoom the sheet to 85% as in the source workbook.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
